# Add a "Username" column (between "No" and "Nama") to the Rekap Hasil Tes
# export sheet, matching the authored OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank column at B (shifts old B..AV to C..AW).
$ws.Columns("B").Insert() | Out-Null

# Give the new column roughly the same width as the header column it now
# sits next to (closest the engine's column-width quantisation allows).
$ws.Columns("B").ColumnWidth = 33.67

# 2. Header label for the new column.
$ws.Range("B7").Value = "Username"

# Copy the header/body formatting from the "Keterangan" column (now E7:E8,
# after the shift) onto the new B7:B8 pair so the new header matches the
# existing look (bordered box, bold, wrap text, etc.).
$ws.Range("E7").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null
$ws.Range("E8").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# Merge the new header cell the same way the other header cells are merged.
$ws.Range("B7:B8").Merge() | Out-Null

# 3. The conditional formatting that greys out empty "Total" values lived on
# column E (E9:E634); after the insert it needs to point at column F.
$cf = $ws.Range("E9:E634").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("F9:F634")) | Out-Null

# 4. Restore the original selected cell shown in the file.
$ws.Range("C13").Select() | Out-Null
